$d = $word.ActiveDocument

# Row-label cells in the financial results table: update text (German label /
# unit tweaks) and turn the row-label run bold (it was w:b w:val="0" i.e.
# explicitly non-bold, now it should be simply <w:b/>).
$replacements = @(
    @{ Old = "Umsatzerlös (k`$)";       New = "Umsatzerlös (`$K)" },
    @{ Old = "Lagerverbrauch (k`$)";    New = "Wareneinsatz (`$K)" },
    @{ Old = "Bruttogewinn (%)";        New = "Bruttogewinnspanne (%)" },
    @{ Old = "Betriebskosten (k`$)";    New = "Betriebskosten (`$K)" },
    @{ Old = "EBITDA (k`$)";            New = "EBITDA (`$K)" },
    @{ Old = "Zinsaufwendungen (k`$)";  New = "Zinsaufwand (`$K)" },
    @{ Old = "Gewinn vor Steuern (k`$)"; New = "Ergebnis vor Steuern (`$K)" },
    @{ Old = "Nettoeinkommen (k`$)";    New = "Nettoeinnahmen (`$K)" },
    @{ Old = "Eigenkapital (k`$)";      New = "Eigenkapital (`$K)" }
)

foreach ($item in $replacements) {
    $d.Content.Find.Execute($item.Old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $item.New, 2) | Out-Null

    $rng = $d.Content.Duplicate
    $rng.Find.Execute($item.New, $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    $rng.Font.Bold = 1
}

# Two more row-label cells whose text stays the same but whose run must
# become bold as well.
$boldOnly = @("Gesamtvermögen (`$K)", "Gesamtverbindlichkeiten (`$K)")

foreach ($text in $boldOnly) {
    $rng = $d.Content.Duplicate
    $rng.Find.Execute($text, $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    $rng.Font.Bold = 1
}
